$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.778.35"
$ws.Range("E2").Value = "'  +6.28%  "
$ws.Range("D3").Value = "'2.303.50"
$ws.Range("E3").Value = "'  +3.37%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'305.94"
$ws.Range("E5").Value = "'  +2.11%  "
$ws.Range("D6").Value = "'102.07"
$ws.Range("E6").Value = "'  +13.14%  "
$ws.Range("E7").Value = "'  +2.67%  "
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "'  +7.38%  "
$ws.Range("D10").Value = "'37.33"
$ws.Range("E10").Value = "'  +13.58%  "
$ws.Range("D11").Value = "'0.0803"
$ws.Range("E11").Value = "'  +3.10%  "
$ws.Range("D12").Value = "'7.46"
$ws.Range("E12").Value = "'  +6.93%  "
$ws.Range("E13").Value = "'  +0.29%  "
$ws.Range("D14").Value = "'2.654.05"
$ws.Range("E14").Value = "'  +3.36%  "
$ws.Range("D15").Value = "'2.305.01"
$ws.Range("D16").Value = "'14.01"
$ws.Range("E16").Value = "'  +3.46%  "
$ws.Range("E17").Value = "'  +5.40%  "
$ws.Range("D18").Value = "'46.768.62"
$ws.Range("E18").Value = "'  +6.46%  "
$ws.Range("D19").Value = "'13.52"
$ws.Range("E19").Value = "'  +21.10%  "
$ws.Range("D20").Value = "'0.0₃0949"
$ws.Range("E20").Value = "'  +4.98%  "
$ws.Range("D21").Value = "'6.10"
$ws.Range("E21").Value = "'  +2.73%  "
$ws.Range("D22").Value = "'66.88"
$ws.Range("E22").Value = "'  +3.29%  "
$ws.Range("D23").Value = "'249.55"
$ws.Range("E23").Value = "'  +4.80%  "
$ws.Range("E24").Value = "'  +4.44%  "
$ws.Range("E25").Value = "'  +5.24%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "'  -0.95%  "
$ws.Range("D27").Value = "'44.07"
$ws.Range("E27").Value = "'  +14.47%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "'  +2.86%  "
$ws.Range("D29").Value = "'9.95"
$ws.Range("E29").Value = "'  +6.55%  "
$ws.Range("D30").Value = "'20.20"
$ws.Range("E30").Value = "'  +4.59%  "
$ws.Range("D31").Value = "'2.87"
$ws.Range("E31").Value = "'  +14.78%  "
$ws.Range("D32").Value = "'5.79"
$ws.Range("E32").Value = "'  +6.96%  "
$ws.Range("B33").Value = "'Monero"
$ws.Range("C33").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'147.53"
$ws.Range("E33").Value = "'  -1.94%  "
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0805"
$ws.Range("E34").Value = "'  +7.25%  "
$ws.Range("D35").Value = "'3.19"
$ws.Range("E35").Value = "'  +12.52%  "
$ws.Range("E36").Value = "'  +11.96%  "
$ws.Range("E37").Value = "'  +2.99%  "
$ws.Range("D38").Value = "'1.81"
$ws.Range("E38").Value = "'  +6.88%  "
$ws.Range("E39").Value = "'  +23.77%  "
$ws.Range("D40").Value = "'4.12"
$ws.Range("E40").Value = "'  +15.10%  "
$ws.Range("D41").Value = "'3.46"
$ws.Range("E41").Value = "'  +7.96%  "
$ws.Range("D42").Value = "'0.0306"
$ws.Range("E42").Value = "'  +1.07%  "
$ws.Range("D43").Value = "'2.03"
$ws.Range("E43").Value = "'  +12.61%  "
$ws.Range("E44").Value = "'  -0.02%  "
$ws.Range("D45").Value = "'1.853.52"
$ws.Range("E45").Value = "'  +1.85%  "
$ws.Range("D46").Value = "'89.08"
$ws.Range("E46").Value = "'  +20.90%  "
$ws.Range("B47").Value = "'ordi"
$ws.Range("C47").Value = "'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "'75.24"
$ws.Range("E47").Value = "'  +12.26%  "
$ws.Range("B48").Value = "'Algorand"
$ws.Range("C48").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.197"
$ws.Range("E48").Value = "'  +9.58%  "
$ws.Range("E49").Value = "'  +10.57%  "
$ws.Range("D50").Value = "'97.14"
$ws.Range("E50").Value = "'  +3.10%  "
$ws.Range("D51").Value = "'54.29"
$ws.Range("E51").Value = "'  +6.69%  "
